# Remove the extra "transaction_status_blockchain" test column (L) that was
# added for the "unexpected column" test fixture, along with its data
# ("test_text_link111" in row 2). Deleting the entire column shifts the
# remaining "Test" value (previously in column M) left into column L,
# matching the target layout (header/data now span A:L instead of A:M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L").Delete()
